$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2248062015503876
$ws.Range("C2").Value = 0.4883720930232558
$ws.Range("J2").Value = 0.02842377260981912
$ws.Range("P2").Value = 0.1550387596899225
$ws.Range("S2").Value = 0.103359173126615
$ws.Range("B3").Value = 0.015625
$ws.Range("C3").Value = 0.02604166666666667
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("O3").Value = 0.005208333333333333
$ws.Range("P3").Value = 0.6770833333333334
$ws.Range("S3").Value = 0.234375
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.1048034934497817
$ws.Range("D6").Value = 0.01746724890829694
$ws.Range("F6").Value = 0.0611353711790393
$ws.Range("J6").Value = 0.2794759825327511
$ws.Range("O6").Value = 0.01746724890829694
$ws.Range("Q6").Value = 0.148471615720524
$ws.Range("R6").Value = 0.07423580786026202
$ws.Range("S6").Value = 0.2969432314410481
$ws.Range("B7").Value = 0.1068376068376068
$ws.Range("D7").Value = 0.0170940170940171
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.141025641025641
$ws.Range("O7").Value = 0.02136752136752137
$ws.Range("Q7").Value = 0.188034188034188
$ws.Range("R7").Value = 0.04273504273504274
$ws.Range("S7").Value = 0.4273504273504273
$ws.Range("B8").Value = 0.109181141439206
$ws.Range("D8").Value = 0.007444168734491315
$ws.Range("F8").Value = 0.06699751861042183
$ws.Range("J8").Value = 0.1513647642679901
$ws.Range("O8").Value = 0.01240694789081886
$ws.Range("Q8").Value = 0.1662531017369727
$ws.Range("R8").Value = 0.0620347394540943
$ws.Range("S8").Value = 0.424317617866005
$ws.Range("B9").Value = 0.128
$ws.Range("D9").Value = 0.028
$ws.Range("E9").Value = 0.004
$ws.Range("F9").Value = 0.064
$ws.Range("J9").Value = 0.152
$ws.Range("O9").Value = 0.012
$ws.Range("Q9").Value = 0.196
$ws.Range("R9").Value = 0.044
$ws.Range("S9").Value = 0.372
$ws.Range("B10").Value = 0.1198581560283688
$ws.Range("D10").Value = 0.02198581560283688
$ws.Range("E10").Value = 0.002127659574468085
$ws.Range("F10").Value = 0.05886524822695036
$ws.Range("J10").Value = 0.1439716312056738
$ws.Range("O10").Value = 0.01276595744680851
$ws.Range("Q10").Value = 0.2304964539007092
$ws.Range("R10").Value = 0.0524822695035461
$ws.Range("S10").Value = 0.3574468085106383
$ws.Range("G11").Value = 0.1457725947521866
$ws.Range("J11").Value = 0.05830903790087463
$ws.Range("K11").Value = 0.2040816326530612
$ws.Range("L11").Value = 0.5743440233236151
$ws.Range("S11").Value = 0.01749271137026239
$ws.Range("G12").Value = 0.7342995169082126
$ws.Range("J12").Value = 0.1739130434782609
$ws.Range("K12").Value = 0.01932367149758454
$ws.Range("L12").Value = 0.03864734299516908
$ws.Range("S12").Value = 0.03381642512077294
$ws.Range("G13").Value = 0.6956521739130435
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.05063291139240506
$ws.Range("H15").Value = 0.1308016877637131
$ws.Range("I15").Value = 0.08016877637130802
$ws.Range("J15").Value = 0.3417721518987342
$ws.Range("K15").Value = 0.09282700421940929
$ws.Range("M15").Value = 0.01687763713080169
$ws.Range("N15").Value = 0.004219409282700422
$ws.Range("O15").Value = 0.04219409282700422
$ws.Range("S15").Value = 0.2405063291139241
$ws.Range("F16").Value = 0.01357466063348416
$ws.Range("H16").Value = 0.1719457013574661
$ws.Range("I16").Value = 0.09954751131221719
$ws.Range("J16").Value = 0.3936651583710407
$ws.Range("K16").Value = 0.09049773755656108
$ws.Range("M16").Value = 0.01809954751131222
$ws.Range("N16").Value = 0.004524886877828055
$ws.Range("O16").Value = 0.05429864253393665
$ws.Range("S16").Value = 0.1538461538461539
$ws.Range("F17").Value = 0.01747572815533981
$ws.Range("H17").Value = 0.1592233009708738
$ws.Range("I17").Value = 0.1281553398058252
$ws.Range("J17").Value = 0.4368932038834951
$ws.Range("K17").Value = 0.08349514563106795
$ws.Range("M17").Value = 0.01941747572815534
$ws.Range("N17").Value = 0.001941747572815534
$ws.Range("O17").Value = 0.06213592233009708
$ws.Range("S17").Value = 0.0912621359223301
$ws.Range("F18").Value = 0.0145985401459854
$ws.Range("H18").Value = 0.1021897810218978
$ws.Range("I18").Value = 0.1094890510948905
$ws.Range("J18").Value = 0.4306569343065693
$ws.Range("K18").Value = 0.1021897810218978
$ws.Range("M18").Value = 0.0145985401459854
$ws.Range("O18").Value = 0.1167883211678832
$ws.Range("S18").Value = 0.1094890510948905
$ws.Range("F19").Value = 0.01783060921248143
$ws.Range("H19").Value = 0.1768202080237742
$ws.Range("I19").Value = 0.09658246656760773
$ws.Range("J19").Value = 0.363298662704309
$ws.Range("K19").Value = 0.1263001485884101
$ws.Range("M19").Value = 0.02005943536404161
$ws.Range("N19").Value = 0.002971768202080238
$ws.Range("O19").Value = 0.07355126300148589
$ws.Range("S19").Value = 0.1225854383358098
